$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Suffrage)
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 0.4878181337611489
$ws.Range("D2").Value = 0.1515156449311006
$ws.Range("E2").Value = "qa_coverage_line_%"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "4.02e-10"

# Row 3 (GatesS)
$ws.Range("C3").Value = 0.4177758339731718
$ws.Range("D3").Value = 0.1334384849465055
$ws.Range("E3").Value = "qa_fix_dispersion_mean"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "1.38e-07"

# Row 4 (GatesT)
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 0.5003541793228794
$ws.Range("D4").Value = 0.1340760916137876
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "7.93e-06"
